$d = $word.ActiveDocument

# Add the missed name "Milad" to the attendance list.
# Original text: ": Dane (Minute Taker), Hamilton (Team Leader), Chris, Luis, Juan"
# New text:      ": Dane (Minute Taker), Hamilton (Team Leader), Chris, Luis, Juan, Milad "
$findRange = $d.Content
$findRange.Find.Execute("Chris, Luis, Juan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertionPoint = $findRange.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.Text = ", Milad "
$newFont = $insertionPoint.Font
$newFont.Name = "Times New Roman"
$newFont.NameBi = "Times New Roman"
$newFont.Size = 14
$newFont.SizeBi = 14
